# Applies the JacobOtaResume.docx edit described by the commit:
#   1. Split "Software Engineering " into "Software Engineering" +
#      " or Software Development Engineer in Test" + " ", leaving the
#      rest of the sentence ("Internships or Entry-Level positions.")
#      as its own separate runs.
#   2. Append ", Information Systems Engineering" and ", Machine Learning"
#      after the "Relevant Coursework" line that ends in "Project Management".
#   3. Append ", SQL" after "...C, C#, Python".
#   4. Append ", MySQL" after "...Git/GitHub, Windows".
#   5. Append ", MySQL Workbench" after "...Visual Studio Code, IntelliJ IDEA".
#
# NOTE on technique: this COM-interop runtime always re-serializes text it
# touches as a single run stretching from the edit point to wherever the
# edit "settles" (end of paragraph for an insert, or the nearer run
# boundary for a delete) UNLESS the touched span is exactly aligned to
# existing run boundaries on both ends. So every edit below either
# (a) inserts brand-new text with the insertion point collapsed to the
#     paragraph's end (the one place an insert reliably becomes its own new
#     <w:r> without disturbing anything earlier in the paragraph), or
# (b) deletes a whole run (run-boundary-to-run-boundary) and then rebuilds
#     the deleted text plus the new text as a sequence of paragraph-end
#     appends, which keeps every neighboring run untouched.

$d = $word.ActiveDocument

# Appends `$newText` as a brand-new run at the very end of `$para`,
# returning nothing. Word COM-interop here only keeps inserted text as a
# separate <w:r> when the insertion point sits at the paragraph's end, so
# we always collapse there first. `$color` (a Long RGB value, e.g. 0 for
# black) is only applied when supplied, since some paragraphs in this
# resume don't stamp an explicit run color.
function Append-RunToParagraph {
    param(
        $para,
        [string]$newText,
        $color
    )

    $r = $para.Range
    $r.Collapse(0)
    $r.MoveEnd(1, -1)        # step back before the paragraph mark
    $r.Collapse(0)
    $insertStart = $r.End

    $r.InsertAfter($newText)

    $newRun = $d.Range($insertStart, $insertStart + $newText.Length)
    $newRun.Font.Size = 10.5
    if ($null -ne $color) {
        $newRun.Font.Color = $color
    }
}

# Finds the paragraph containing `$findText` and appends `$newText` (a new
# run) right after it, at the paragraph's end.
function Append-RunAtParagraphEnd {
    param(
        [string]$findText,
        [string]$newText,
        $color
    )

    $found = $d.Content
    $null = $found.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $para = $found.Paragraphs(1)
    Append-RunToParagraph $para $newText $color
    return $para
}

# --- Edit 1: Summary line -------------------------------------------------
# "...seeking Software Engineering Internships or Entry-Level positions."
# becomes
# "...seeking Software Engineering or Software Development Engineer in Test
#  Internships or Entry-Level positions."
#
# The whole "Software Engineering Internships or Entry-Level positions."
# span is exactly the concatenation of five whole runs, so deleting it is a
# clean run-boundary-to-run-boundary delete that leaves the preceding
# "...seeking " run completely untouched. We then rebuild it plus the new
# wording as a sequence of paragraph-end appends (each becomes its own
# <w:r>).

$summaryMatch = $d.Content
$null = $summaryMatch.Find.Execute("Software Engineering Internships or Entry-Level positions.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$summaryPara = $summaryMatch.Paragraphs(1)

$summaryDelete = $d.Range($summaryMatch.Start, $summaryMatch.End)
$summaryDelete.Delete()

Append-RunToParagraph $summaryPara "Software Engineering" 0
Append-RunToParagraph $summaryPara " or Software Development Engineer in Test" 0
Append-RunToParagraph $summaryPara " " 0
Append-RunToParagraph $summaryPara "Internship" 0
Append-RunToParagraph $summaryPara "s" 0
Append-RunToParagraph $summaryPara " " 0
Append-RunToParagraph $summaryPara "or Entry-Level positions" 0
Append-RunToParagraph $summaryPara "." 0

# --- Edit 2: Relevant Coursework (Masters) --------------------------------
$null = Append-RunAtParagraphEnd "Software Testing, Distributed Software Development and Integration, Project Management" ", Information Systems Engineering" $null
$null = Append-RunAtParagraphEnd "Software Testing, Distributed Software Development and Integration, Project Management" ", Machine Learning" $null

# --- Edit 3: Programming Languages ----------------------------------------
$null = Append-RunAtParagraphEnd ", Python" ", SQL" 0

# --- Edit 4: Tools, Frameworks, and OS ------------------------------------
$null = Append-RunAtParagraphEnd "GitHub, Windows" ", MySQL" 0

# --- Edit 5: IDEs ----------------------------------------------------------
$null = Append-RunAtParagraphEnd "Visual Studio, Visual Studio Code, IntelliJ IDEA" ", MySQL Workbench" 0

Write-Output "Done."
